$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rename the sheet
$ws.Name = "participantes"

# Remove the stray "applyBorder" formatting (no visible border) from column C.
$ws.Range("C2:C38").Borders.LineStyle = -4142

# Same cleanup on column D, but only for the rows that don't carry the
# hyperlink style - those keep their font/style untouched.
$ws.Range("D2").Borders.LineStyle = -4142
$ws.Range("D4:D11").Borders.LineStyle = -4142
$ws.Range("D14:D16").Borders.LineStyle = -4142
$ws.Range("D18:D24").Borders.LineStyle = -4142
$ws.Range("D26:D28").Borders.LineStyle = -4142
$ws.Range("D30").Borders.LineStyle = -4142
$ws.Range("D32").Borders.LineStyle = -4142
$ws.Range("D34").Borders.LineStyle = -4142
$ws.Range("D36:D37").Borders.LineStyle = -4142

# Clear out the leftover helper columns L:M
$ws.Range("L1:M38").Clear()

# Remove the now-empty trailing row
$ws.Rows("38:38").Delete()

# Match the saved selection/active cell
$ws.Range("C42").Select()
